$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells we will touch, so numeric-looking
# strings (e.g. "1.00", "7.29") are preserved exactly as text like the original file.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.685.54"
$ws.Range("E2").Value = "  +1.84%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.935.72"
$ws.Range("E3").Value = "  -0.12%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.03"
$ws.Range("E5").Value = "  -0.82%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.19"
$ws.Range("E6").Value = "  +0.56%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.933.20"
$ws.Range("E8").Value = "  -0.17%  "

# Row 9
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.505"
$ws.Range("E9").Value = "  +0.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").Value = "  +4.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  +5.13%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("E12").Value = "  +0.08%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000236"
$ws.Range("E13").Value = "  +4.57%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.58"
$ws.Range("E14").Value = "  -3.21%  "

# Row 15
$ws.Range("E15").Value = "  -0.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.423.50"
$ws.Range("E16").Value = "  -0.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.675.47"
$ws.Range("E17").Value = "  +1.80%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.66"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.938.71"
$ws.Range("E19").Value = "  +0.06%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "438.14"
$ws.Range("E20").Value = "  +1.26%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.32"
$ws.Range("E21").Value = "  -1.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.663"
$ws.Range("E22").Value = "  -2.34%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.01"
$ws.Range("E23").Value = "  -1.43%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.87"
$ws.Range("E24").Value = "  -1.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.00"
$ws.Range("E25").Value = "  +0.59%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.12"
$ws.Range("E26").Value = "  -3.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.69"
$ws.Range("E27").Value = "  -1.07%  "

# Row 28
$ws.Range("E28").Value = "  -0.05%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  -0.23%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.17"
$ws.Range("E30").Value = "  +3.15%  "

# Row 31
$ws.Range("E31").Value = "  -0.51%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000101"
$ws.Range("E32").Value = "  +14.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.108"
$ws.Range("E33").Value = "  -1.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.30"
$ws.Range("E34").Value = "  -1.53%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.993"
$ws.Range("E36").Value = "  -1.93%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.09"
$ws.Range("E37").Value = "  +2.70%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.54"
$ws.Range("E38").Value = "  -1.89%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.61"
$ws.Range("E39").Value = "  -0.88%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.01"
$ws.Range("E40").Value = "  +0.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.44"
$ws.Range("E41").Value = "  -1.79%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.279"
$ws.Range("E43").Value = "  -1.51%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.98"
$ws.Range("E44").Value = "  -8.03%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.696.43"
$ws.Range("E45").Value = "  -0.21%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.31"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "360.46"
$ws.Range("E47").Value = "  -1.25%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0334"
$ws.Range("E48").Value = "  -3.74%  "

# Row 50
$ws.Range("E50").Value = "  -1.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.65"
$ws.Range("E51").Value = "  -4.84%  "

Write-Host "Applied crypto price updates"